$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 8

# Copy formatting (incl. number format/style) from the row above so the
# new row's date cell gets the same style index as existing date cells,
# without Excel fabricating a brand-new custom number format entry.
$ws.Range("A" + ($row - 1)).Copy()
$ws.Range("A" + $row).PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item($row, 1).Value = 42608.888761574075
$ws.Cells.Item($row, 2).Value = 2
$ws.Cells.Item($row, 3).Value = 48
$ws.Cells.Item($row, 4).Value = 42
$ws.Cells.Item($row, 5).Value = 51
$ws.Cells.Item($row, 6).Value = 48
$ws.Cells.Item($row, 7).Value = 15674
$ws.Cells.Item($row, 8).Value = 32410
$ws.Cells.Item($row, 9).Value = 3674
$ws.Cells.Item($row, 10).Value = 487
$ws.Cells.Item($row, 11).Value = 427
$ws.Cells.Item($row, 12).Value = 14
$ws.Cells.Item($row, 13).Value = 13
$ws.Cells.Item($row, 14).Value = "Bag"
